$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'60.693.00"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -4.12%  "

$ws.Cells.Item(3, 4).Value = "'2.914.51"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -3.72%  "

$ws.Cells.Item(4, 5).Value = "  +0.07%  "

$ws.Cells.Item(5, 4).Value = "'526.58"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -5.86%  "

$ws.Cells.Item(6, 4).Value = "'143.83"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -7.83%  "

$ws.Cells.Item(7, 5).Value = "  +0.12%  "

$ws.Cells.Item(8, 4).Value = "'0.555"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.52%  "

$ws.Cells.Item(9, 4).Value = "'2.921.97"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.80%  "

$ws.Cells.Item(10, 4).Value = "'0.108"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -5.33%  "

$ws.Cells.Item(11, 4).Value = "'5.88"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -8.56%  "

$ws.Cells.Item(12, 4).Value = "'0.353"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -4.07%  "

$ws.Cells.Item(13, 4).Value = "'3.425.07"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -3.58%  "

$ws.Cells.Item(14, 5).Value = "  +1.06%  "

$ws.Cells.Item(15, 4).Value = "'60.793.16"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -3.99%  "

$ws.Cells.Item(16, 4).Value = "'22.73"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -6.14%  "

$ws.Cells.Item(17, 4).Value = "'2.918.11"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.46%  "

$ws.Cells.Item(18, 4).Value = "'0.0000140"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -7.59%  "

$ws.Cells.Item(19, 4).Value = "'4.96"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.08%  "

$ws.Cells.Item(20, 4).Value = "'11.61"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.90%  "

$ws.Cells.Item(21, 4).Value = "'360.90"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -9.67%  "

$ws.Cells.Item(22, 4).Value = "'6.48"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.30%  "

$ws.Cells.Item(23, 4).Value = "'0.999"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.15%  "

$ws.Cells.Item(24, 4).Value = "'5.66"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.47%  "

$ws.Cells.Item(25, 4).Value = "'64.15"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.00%  "

$ws.Cells.Item(26, 4).Value = "'3.047.70"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -3.51%  "

$ws.Cells.Item(27, 4).Value = "'0.451"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -3.59%  "

$ws.Cells.Item(28, 4).Value = "'0.183"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.96%  "

$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.23%  "

$ws.Cells.Item(30, 4).Value = "'0.0₃0861"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -13.38%  "

$ws.Cells.Item(31, 4).Value = "'7.65"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -12.58%  "

$ws.Cells.Item(32, 5).Value = "  -0.02%  "

$ws.Cells.Item(33, 4).Value = "'1.66"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.59%  "

$ws.Cells.Item(34, 4).Value = "'19.70"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -3.92%  "

$ws.Cells.Item(35, 4).Value = "'154.44"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -4.25%  "

$ws.Cells.Item(36, 4).Value = "'4.35"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -8.92%  "

$ws.Cells.Item(37, 4).Value = "'5.61"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -7.91%  "

$ws.Cells.Item(38, 4).Value = "'1.00"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -10.56%  "

$ws.Cells.Item(39, 4).Value = "'1.21"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -8.51%  "

$ws.Cells.Item(40, 4).Value = "'38.05"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.90%  "

$ws.Cells.Item(41, 2).Value = "Maker"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(41, 4).Value = "'2.341.26"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -8.26%  "

$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).Value = "'1.47"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -7.71%  "

$ws.Cells.Item(43, 4).Value = "'3.70"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -6.83%  "

$ws.Cells.Item(44, 4).Value = "'0.647"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.89%  "

$ws.Cells.Item(45, 4).Value = "'20.71"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -9.98%  "

$ws.Cells.Item(46, 4).Value = "'0.0568"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -5.86%  "

$ws.Cells.Item(47, 5).Value = "  +0.20%  "

$ws.Cells.Item(48, 4).Value = "'4.86"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -4.80%  "

$ws.Cells.Item(49, 4).Value = "'0.0233"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -7.29%  "

$ws.Cells.Item(50, 4).Value = "'10.35"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.41%  "

$ws.Cells.Item(51, 4).Value = "'0.0924"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.63%  "
